# "Loan RBI, Variable Instalments"
#
# The "Repayment Schedule" sheet gains a new (blank) column at N, pushing
# the previous N/O/P ("In Advance"(M stays) / "Late" / "Heading" / "Outstanding")
# columns one slot to the right, and that sheet becomes the active tab with
# the cursor resting on U11 (scrolled so column F is near the left edge).
# The previously-active "NewLoanInput" sheet is no longer the selected tab.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before column N - this shifts the old
# N ("Late"), O ("Heading"), P ("Outstanding") columns to O, P, Q.
$ws.Columns("N").Insert()

# Make "Repayment Schedule" the active sheet/tab.
$ws.Activate()

# Scroll so column F is the left-most visible column (topLeftCell = F1),
# then rest the selection on U11.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 6
$ws.Range("U11").Select()
